# Sprint_1_folders/Feasibility_Study.docx
# "updating the study based on suggestions from sponsor"
#
# The final paragraph of the document currently reads ".  " (a period
# followed by the two non-breaking spaces that were already there). The
# sponsor update replaces that period with a new sentence about moving to
# a mobile application, while the original trailing whitespace is kept
# (now as its own run), and "bluetooth" gets wrapped in spell-check
# proofErr markers, just like Word does when it flags a word as a
# possible misspelling.

$d = $word.ActiveDocument

# Locate the paragraph whose whole text is just ".  " (period + two
# non-breaking spaces) -- this is the paragraph that needs the update.
$target = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]0x0D, [char]0x07)
    if ($t -eq ("." + [char]0x00A0 + [char]0x00A0)) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph (period + nbsp nbsp)."
}

# Range covering just the paragraph's own content, excluding the
# paragraph mark at the end, so the paragraph itself (and its pPr/rsid
# attributes) is left untouched -- only its runs are replaced. Rebuilt
# via Document.Range(start, end) (rather than reusing the paragraph's
# own Range object directly) since InsertXML needs a plain start/end
# range to replace in place.
$full = $target.Range
$body = $d.Range($full.Start, $full.End - 1)

$nbsp = [char]0x00A0
$tail = "$nbsp$nbsp"

$innerXml = '<w:r><w:t xml:space="preserve">Update: Based on our meeting with our sponsor, we have decided to move from a web based application, to a mobile application. This does not affect the estimation given above as the technology stack we will be using will still remain open sourced. The budget will also remain relatively low as there are cheap </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>bluetooth</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> beacons available that we can utilize. </w:t></w:r>' + `
            '<w:r><w:t>' + $tail + '</w:t></w:r>'

$packageXml = '<?xml version="1.0" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
            '<w:body><w:p>' + $innerXml + '</w:p></w:body>' + `
          '</w:document>' + `
        '</pkg:xmlData>' + `
      '</pkg:part>' + `
    '</pkg:package>'

$body.InsertXML($packageXml)
